$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44433
$ws.Range("M2").Value2 = 15

$ws.Range("D3").Value2 = 44377
$ws.Range("M3").Value2 = 30
$ws.Range("N3").Value2 = 40000
$ws.Range("O3").Value2 = 40000
$ws.Range("P3").Value2 = 40000
$ws.Range("S3").Value2 = 2222

$ws.Range("D4").Value2 = 44435
$ws.Range("R4").Value2 = 'Perú'

$ws.Range("D5").Value2 = 44435
$ws.Range("M5").Value2 = 105
$ws.Range("N5").Value2 = 35000
$ws.Range("O5").Value2 = 35000
$ws.Range("P5").Value2 = 35000
$ws.Range("S5").Value2 = 1944

$ws.Range("D6").Value2 = 44438
$ws.Range("M6").Value2 = 25
$ws.Range("O6").Value2 = 35000
$ws.Range("P6").Value2 = 35000
$ws.Range("S6").Value2 = 1944

$ws.Range("D8").Value2 = 44369
$ws.Range("M8").Value2 = 5
$ws.Range("R8").Value2 = 'Perú'

$ws.Range("D9").Value2 = 44442
$ws.Range("M9").Value2 = 15
$ws.Range("N9").Value2 = 35000
$ws.Range("O9").Value2 = 35000
$ws.Range("P9").Value2 = 35000
$ws.Range("Q9").Value2 = '$/caja 18 kilos'
$ws.Range("R9").Value2 = 'Perú'
$ws.Range("S9").Value2 = 1944
$ws.Range("T9").Value2 = 18

$ws.Range("D10").Value2 = 44363
$ws.Range("M10").Value2 = 144
$ws.Range("N10").Value2 = 1700
$ws.Range("O10").Value2 = 1700
$ws.Range("P10").Value2 = 1700
$ws.Range("Q10").Value2 = '$/kilo'
$ws.Range("R10").Value2 = 'Región de Arica y Parinacota'
$ws.Range("S10").Value2 = 1700
$ws.Range("T10").Value2 = 1

$ws.Range("D11").Value2 = 44392
$ws.Range("M11").Value2 = 20
$ws.Range("R11").Value2 = 'Región de Arica y Parinacota'

$ws.Range("D12").Value2 = 44431
$ws.Range("M12").Value2 = 30
$ws.Range("N12").Value2 = 35000
$ws.Range("O12").Value2 = 35000
$ws.Range("P12").Value2 = 35000
$ws.Range("S12").Value2 = 1944

$ws.Range("D13").Value2 = 44424
$ws.Range("M13").Value2 = 15

$ws.Range("D14").Value2 = 44294
$ws.Range("M14").Value2 = 15
$ws.Range("N14").Value2 = 35000
$ws.Range("O14").Value2 = 35000
$ws.Range("P14").Value2 = 35000
$ws.Range("S14").Value2 = 1944

$ws.Range("D15").Value2 = 44364
$ws.Range("M15").Value2 = 90

$ws.Range("D16").Value2 = 44379
$ws.Range("M16").Value2 = 10
$ws.Range("N16").Value2 = 30000
$ws.Range("O16").Value2 = 30000
$ws.Range("P16").Value2 = 30000
$ws.Range("S16").Value2 = 1667

$ws.Range("D17").Value2 = 44449
$ws.Range("M17").Value2 = 20
$ws.Range("N17").Value2 = 38000
$ws.Range("O17").Value2 = 38000
$ws.Range("P17").Value2 = 38000
$ws.Range("R17").Value2 = 'Región de Arica y Parinacota'
$ws.Range("S17").Value2 = 2111

$ws.Range("D18").Value2 = 44418
$ws.Range("M18").Value2 = 30
$ws.Range("R18").Value2 = 'Región de Arica y Parinacota'

$ws.Range("D19").Value2 = 44434
$ws.Range("M19").Value2 = 40

$ws.Range("D20").Value2 = 44279
$ws.Range("M20").Value2 = 30
$ws.Range("N20").Value2 = 35000
$ws.Range("O20").Value2 = 36000
$ws.Range("P20").Value2 = 35667
$ws.Range("R20").Value2 = 'Región de Arica y Parinacota'
$ws.Range("S20").Value2 = 1982

$ws.Range("D21").Value2 = 44264
$ws.Range("M21").Value2 = 20

$ws.Range("D22").Value2 = 44405
$ws.Range("M22").Value2 = 10

$ws.Range("D23").Value2 = 44432
$ws.Range("M23").Value2 = 10
$ws.Range("R23").Value2 = 'Perú'

$ws.Range("D24").Value2 = 44357
$ws.Range("M24").Value2 = 10
$ws.Range("N24").Value2 = 38000
$ws.Range("O24").Value2 = 38000
$ws.Range("P24").Value2 = 38000
$ws.Range("R24").Value2 = 'Perú'
$ws.Range("S24").Value2 = 2111
